# Updated symbol list (Price / Volume(1h) columns) for the crypto table.
# Values are stored as literal text in the source workbook (t="inlineStr"),
# so each write is prefixed with a leading apostrophe to force Excel to
# keep them as text instead of auto-converting to numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'299.63"
$ws.Range("E2").Value = "'-4.83%"

$ws.Range("D3").Value = "'35.10"
$ws.Range("E3").Value = "'-0.49%"

$ws.Range("D4").Value = "'5.032"
$ws.Range("E4").Value = "'-1.17%"

$ws.Range("E5").Value = "'-2.73%"

$ws.Range("D6").Value = "'1.894"
$ws.Range("E6").Value = "'-9.06%"

$ws.Range("D7").Value = "'7.779"
$ws.Range("E7").Value = "'-2.04%"

$ws.Range("D8").Value = "'4.031"
$ws.Range("E8").Value = "'-2.56%"

$ws.Range("D9").Value = "'0.9257"
$ws.Range("E9").Value = "'-0.61%"

$ws.Range("D10").Value = "'0.1426"
$ws.Range("E10").Value = "'37.81%"

$ws.Range("D11").Value = "'0.1900"
$ws.Range("E11").Value = "'-1.67%"

$ws.Range("D12").Value = "'0.09213"
$ws.Range("E12").Value = "'1.25%"

$ws.Range("D13").Value = "'0.03462"
$ws.Range("E13").Value = "'-3.78%"

$ws.Range("D14").Value = "'0.09867"
$ws.Range("E14").Value = "'-0.27%"

$ws.Range("D15").Value = "'0.001397"
$ws.Range("E15").Value = "'-2.50%"

$ws.Range("D16").Value = "'0.005717"
$ws.Range("E16").Value = "'0.72%"

$ws.Range("D17").Value = "'3.519"
$ws.Range("E17").Value = "'1.44%"

$ws.Range("D18").Value = "'2.962"
$ws.Range("E18").Value = "'3.80%"

$ws.Range("D19").Value = "'0.3410"
$ws.Range("E19").Value = "'-1.41%"

$ws.Range("D20").Value = "'0.1292"
$ws.Range("E20").Value = "'-0.54%"

$ws.Range("D21").Value = "'5.039"
$ws.Range("E21").Value = "'-1.23%"

$ws.Range("D23").Value = "'0.04463"
$ws.Range("E23").Value = "'-1.82%"

$ws.Range("D24").Value = "'0.001216"
$ws.Range("E24").Value = "'-2.13%"

$ws.Range("D25").Value = "'0.004754"
$ws.Range("E25").Value = "'-0.91%"

$ws.Range("D26").Value = "'0.0001232"
$ws.Range("E26").Value = "'-1.58%"

$ws.Range("D27").Value = "'0.0003000"
$ws.Range("E27").Value = "'-33.30%"

$ws.Range("D39").Value = "'0.01897"
$ws.Range("E39").Value = "'-4.52%"

$ws.Range("D40").Value = "'0.04711"
$ws.Range("E40").Value = "'-4.90%"

$ws.Range("D41").Value = "'0.007358"
$ws.Range("E41").Value = "'-3.28%"

$ws.Range("D42").Value = "'0.009681"
$ws.Range("E42").Value = "'23.02%"

$ws.Range("D43").Value = "'0.1316"
$ws.Range("E43").Value = "'-5.02%"

$ws.Range("D44").Value = "'0.002113"
$ws.Range("E44").Value = "'-6.61%"

$ws.Range("D45").Value = "'0.009327"
$ws.Range("E45").Value = "'-20.66%"

$ws.Range("D46").Value = "'0.00006239"
$ws.Range("E46").Value = "'-6.53%"

$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.13%"

$ws.Range("D48").Value = "'64.53"
$ws.Range("E48").Value = "'-65.73%"

$ws.Range("D49").Value = "'0.001659"
$ws.Range("E49").Value = "'-2.39%"

$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.13%"

$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.13%"
